# Laborator 09.11.2023: Liste create de utilizator si parcurgerile in graf
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the new student "Pintilie Robert" into the first empty row (row 43),
#    and mark his attendance for "sapt 6" (column H).
$ws.Range("B43").Value = "Pintilie Robert"
$ws.Range("H43").Value = $true

# 2) Mark attendance ("sapt 6", column H) for several existing students.
$ws.Range("H7").Value = $true
$ws.Range("H8").Value = $true
$ws.Range("H13").Value = $true
$ws.Range("H14").Value = $true
$ws.Range("H25").Value = $true
$ws.Range("H27").Value = $true
$ws.Range("H28").Value = $true
$ws.Range("H29").Value = $true

# 3) Re-sort the roster (B3:S43) alphabetically by name (column B), which puts
#    "Pintilie Robert" into its alphabetical slot and pushes later rows down.
$sortRange = $ws.Range("B3:S43")
$key1 = $ws.Range("B43")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# 4) Update the view: scroll so row 19 is at the top and select J38.
$ws.Range("J38").Select()
$excel.ActiveWindow.ScrollRow = 19
